# "corrected output board bom"
# R10..R17 (rows 20-27) had mistakenly unique per-part Device/Package values
# (R-EU_M0806/M0806 .. R-EU_M0813/M0813). Correct them back to the same
# Device/Package used by the other 0805 resistors (R-EU_M0805/M0805),
# matching rows such as R1..R9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 20; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = "R-EU_M0805"
    $ws.Cells.Item($row, 4).Value = "M0805"
}

# Update the last-selected cell to match the saved workbook state.
$ws.Range("F10").Select()
